$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "B11" "68.28"
Set-TextValue "C11" "3.48"
Set-TextValue "D11" "71.75"

Set-TextValue "C12" "33.38"
Set-TextValue "D12" "80.58"

Set-TextValue "B33" "59.98"
Set-TextValue "C33" "3.27"
Set-TextValue "D33" "63.24"

Set-TextValue "B34" "45.56"
Set-TextValue "C34" "34.17"
Set-TextValue "D34" "79.74"

Set-TextValue "B36" "94.76"
Set-TextValue "C36" "5.16"
Set-TextValue "D36" "99.92"

Set-TextValue "C40" "38.35"
Set-TextValue "D40" "67.95"
